# Insert 3 new rows of weekly price data right before the current row 641
# (Fecha serial 45265 = 2023-12-05), shifting the existing rows 641-694 down
# to 644-697.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 641..643
$ws.Range("A641:A643").EntireRow.Insert()

# Common values shared across these rows
$mercadoId = 11
$mercado   = "Vega Monumental Concepción"
$region    = "Bíobío"
$codreg    = 8
$catId     = 100112006
$categoria = "Repollo"
$unidad    = "$/unidad"
$origen    = "Región Metropolitana"
$kgUnid    = 1
$clasif    = "Hortaliza"
$fecha     = 45265

function Set-Fila($fila, $variedad, $calidad, $volumen, $precioMin, $precioMax, $precioProm) {
    $ws.Cells.Item($fila, 1).Value  = $mercadoId
    $ws.Cells.Item($fila, 2).Value  = $mercado
    $ws.Cells.Item($fila, 3).Value  = $region
    $ws.Cells.Item($fila, 4).Value  = $fecha
    $ws.Cells.Item($fila, 5).Value  = $codreg
    $ws.Cells.Item($fila, 6).Value  = $catId
    $ws.Cells.Item($fila, 7).Value  = $categoria
    $ws.Cells.Item($fila, 8).Value  = $variedad
    $ws.Cells.Item($fila, 9).Value  = $calidad
    $ws.Cells.Item($fila, 10).Value = $volumen
    $ws.Cells.Item($fila, 11).Value = $precioMin
    $ws.Cells.Item($fila, 12).Value = $precioMax
    $ws.Cells.Item($fila, 13).Value = $precioProm
    $ws.Cells.Item($fila, 14).Value = $unidad
    $ws.Cells.Item($fila, 15).Value = $origen
    $ws.Cells.Item($fila, 16).Value = $precioProm
    $ws.Cells.Item($fila, 17).Value = $kgUnid
    $ws.Cells.Item($fila, 18).Value = $clasif
}

Set-Fila 641 "Crespo record" "Primera" 1000 900  1000 950
Set-Fila 642 "Crespo record" "Segunda" 500  800  800  800
Set-Fila 643 "Morada(o)"     "Primera" 500  1200 1200 1200
